$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that used to sit right
#    under the title heading (2nd paragraph of the document).
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Almighty Reels Power of Zeus Free |
#    Game Review" right before the final paragraph (the old "Prompt for
#    DALLE..." paragraph), without disturbing that paragraph's own run
#    structure/formatting.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$insertPos = $lastPara.Range.Start
$insertPoint = $d.Range($insertPos, $insertPos)

$xmlFragment = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Almighty Reels Power of Zeus Free | Game Review</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$null = $insertPoint.InsertXML($xmlFragment)

# The fragment above carries a trailing empty paragraph mark so that the
# original last paragraph keeps its own formatting/runs untouched; remove
# that now-superfluous empty paragraph. (An "empty" paragraph's Range.Text
# is just the paragraph mark itself, i.e. length 1.)
$emptyPara = $d.Paragraphs($count + 1)
if ($emptyPara.Range.Text.Length -le 1) {
    $emptyPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Replace the text of the old "Prompt for DALLE: ..." paragraph (now the
#    very last paragraph) with the meta-description copy, keeping its
#    existing italic run formatting intact.
# ---------------------------------------------------------------------------
$oldText = 'Prompt for DALLE: Create a cartoon style feature image for the game "Almighty Reels Power of Zeus". The image should feature a happy Maya warrior with glasses.'
$newText = 'Read our review of Almighty Reels Power of Zeus online slot game. Play for free and learn about the gameplay, symbols, bonuses, and similar games.'

$null = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                                 $true, 1, $false, $newText, 2)
